$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Activate()

# --- Update the "stable coin" contract address (shared string) ---
$ws.Range("H2").Value = "0x3c499c542cef5e3811e1192ce70d8cc03d5c3359"

# --- Update the dollar-increment multiplier ---
$ws.Range("D2").Value = 2

# --- Update the first data row (A4:C4): amount bought, running sum seed, date ---
$ws.Range("A4").Value = 9999
$ws.Range("B4").Value = 0.0000000000001
$ws.Range("C4").Value = 45972

# --- Recalculate dependent formulas (I2, J2, K2, L2, M2) ---
$excel.Calculate()

# --- Update the active selection / view state ---
$ws.Range("A5:D22").Select()
